$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores prices/percentages as plain text (inline strings). Excel
# normally "helpfully" reinterprets numeric-looking text as a Double when a
# Range.Value is assigned, which would mangle values like "215.83" ->
# "215.83000000000001" or drop significant formatting like trailing zeros.
# Force the target range to Text format while writing, then restore the
# original (default/"Normal") cell style so no stray formatting is left behind.
$dataRng = $ws.Range("B2:E51")
$dataRng.NumberFormat = "@"

$ws.Range("D2").Value = "25.882.66"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "1.639.63"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "215.83"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "0.5039"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.2566"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "0.06397"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "19.70"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "4.265"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "1.639.15"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "1.864.89"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "0.0₅7909"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "64.36"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "25.918.42"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "203.35"
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("D21").Value = "4.381"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "9.927"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "5.976"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "1.937"
$ws.Range("E25").Value = "  +11.01%  "
$ws.Range("D26").Value = "141.23"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").Value = "0.1138"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").Value = "15.71"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "6.769"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "0.04948"
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "3.189"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "1.548"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "2.381"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "2.629"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").Value = "0.8932"
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("D38").Value = "1.162.44"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "0.5619"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "0.01570"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "1.007"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "5.661"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "100.00"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8085"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").Value = "1.776.22"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").Value = "0.4530"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "54.97"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "0.05054"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  -0.04%  "

$dataRng.Style = "Normal"
